# Remove the "Justification (à l'oral ou en note)" column (the 6th /
# last column) from the scenario-probability table on slide 6.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTable) {
        $tbl = $shp.Table
        $tbl.Columns.Item($tbl.Columns.Count).Delete()
    }
}
